# desgn_params.xlsx update: add report-table "Sheet2" (Operations / Parameters),
# update orbiter_mass value on PRIMARY INPUTS, and touch a few selection/format
# details to match the authored commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. PRIMARY INPUTS: update orbiter_mass (B2) 3284.166976 -> 5438.66
# ---------------------------------------------------------------------------
$primary = $wb.Worksheets.Item("PRIMARY INPUTS")
$primary.Range("B2").Value = 5438.66

# Column width tweaks on PRIMARY INPUTS (col A wider, new col E width)
$primary.Columns.Item(1).ColumnWidth = 17.6
$primary.Columns.Item(5).ColumnWidth = 22.7

# New (empty, bold-formatted) row 13 below the existing data
$primary.Rows.Item(13).Font.Bold = $true

# Move the sheet's active selection the way the authored file shows it
$primary.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert a new worksheet named "Sheet2" right after "PRIMARY INPUTS"
#    holding the "Operations" / "Parameters" report tables.
# ---------------------------------------------------------------------------
$ops = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $primary)
$ops.Name = "Sheet2"

# -- section headers --------------------------------------------------------
$ops.Range("A1").Value = "Operations"
$ops.Range("E1").Value = "Parameters"

# -- table headers (bold) ----------------------------------------------------
$ops.Range("A3:G3").Font.Bold = $true
$ops.Range("A3").Value = "name"
$ops.Range("B3").Value = "value"
$ops.Range("E3").Value = "name"
$ops.Range("F3").Value = "value"

# -- left ("Operations") table: mirrors orb_mission A2:C11 ------------------
$ops.Range("A4").Value = "mom_dump_freq"
$ops.Range("B4").Value = 0.5
$ops.Range("C4").Value = "per day"

$ops.Range("A5").Value = "burn_time"
$ops.Range("B5").Value = 1
$ops.Range("C5").Value = "sec"

$ops.Range("A6").Value = "lifetime"
$ops.Range("B6").Value = 5
$ops.Range("C6").Value = "years"

$ops.Range("A7").Value = "n_slew_maneuvers"
$ops.Range("B7").Value = 5

$ops.Range("A8").Value = "slew_angle"
$ops.Range("B8").Value = 90
$ops.Range("C8").Value = "deg"

$ops.Range("A9").Value = "slew_time"
$ops.Range("B9").Value = 3000
$ops.Range("C9").Value = "sec"

$ops.Range("A10").Value = "slew_burn_pct"
$ops.Range("B10").Value = 0.05
$ops.Range("C10").Value = 0.05
$ops.Range("C10").NumberFormat = "0%"

$ops.Range("A11").Font.Bold = $true
$ops.Range("A11").Value = "app_slew"
$ops.Range("B11").Value = 90
$ops.Range("C11").Value = "deg"

$ops.Range("A12").Font.Bold = $true
$ops.Range("A12").Value = "app_time"
$ops.Range("B12").Value = 120
$ops.Range("C12").Value = "sec"

# -- right ("Parameters") table ---------------------------------------------
$ops.Range("E4").Value = "dipole"
$ops.Range("F4").Value = 1
$ops.Range("G4").Value = "A m^2"

$ops.Range("E5").Value = "solar incidence"
$ops.Range("F5").Value = 0
$ops.Range("G5").Value = "deg"

$ops.Range("E6").Value = "Cd"
$ops.Range("F6").Value = 2.2

$ops.Range("E7").Value = "q"
$ops.Range("F7").Value = 0.6

$ops.Range("E8").Value = "pt accuracy"
$ops.Range("F8").Value = 0.5
$ops.Range("G8").Value = "deg"

$ops.Range("E9").Value = "orbiter_mass"
$ops.Range("F9").Formula = "='PRIMARY INPUTS'!B2"
$ops.Range("G9").Value = "kg"

$ops.Range("E10").Value = "solar array area"
$ops.Range("F10").Value = 7.8826
$ops.Range("G10").Value = "m^2"

$ops.Range("E11").Value = "TTC earth pt array moi"
$ops.Range("F11").Value = 102.87
$ops.Range("G11").Value = "m^4"

# Match the authored selection state for the new sheet
$ops.Range("A1:G13").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Cosmetic selection updates on the other (shifted) worksheets
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("orb_mission").Range("H23").Select() | Out-Null
$wb.Worksheets.Item("orb_props").Range("A1:C8").Select() | Out-Null

# Re-activate PRIMARY INPUTS so it stays the displayed/active tab
$primary.Activate() | Out-Null
$primary.Range("B5").Select() | Out-Null
